$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Split the "names" column into "name" (first name) + "last_name" ---
# Insert a blank column before the existing "names" column (A), shifting
# names/unique_number/email/career_id one column to the right (B..E).
$ws.Columns("A").Insert()

# The inserted column doesn't pick up the old column A width automatically,
# so restore it to match what column B (the old "names" column) has.
$ws.Columns("A").ColumnWidth = $ws.Columns("B").ColumnWidth

# The worksheet's hyperlinks (originally anchored on the "email" column, C)
# stayed anchored on column C after the insert instead of moving to D, so
# drop and re-add them on the correct (shifted) column.
$ws.Range("C2").Hyperlinks.Delete()

# New header row.
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "last_name"

# Row 2 - Jonathan Vasquez
$ws.Range("A2").Value = "Jonathan "
$ws.Range("B2").Value = "Vasquez"

# Fill in last names for the remaining rows first...
$ws.Range("B3").Value = "Morales"
$ws.Range("B4").Value = "Zambrano"
$ws.Range("B5").Value = "Segovia"

# ...then come back and fill in the first names.
$ws.Range("A3").Value = "Chantal"
$ws.Range("A4").Value = "Nicole"
$ws.Range("A5").Value = "Kevin"

# Re-create the mailto hyperlinks against the shifted "email" column (D).
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:jonathan.vasquez01@epn.edu.ec")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:chantal.morales@epn.edu.ec")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:nicole.zambrano@epn.edu.ec")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:kevin.segovia@epn.edu.ec")

$ws.Range("C10").Select()
